# Apply changes to input cost data on the "time_variants" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time_variants")

# Update row 5 (econ_program_totalcost_vaccination) values
$ws.Range("J5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("Q5").Value = 30000
$ws.Range("V5").Value = 50000
$ws.Range("AF5").Value = 80000
$ws.Range("AK5").Value = 100000
$ws.Range("AP5").Value = 100000
$ws.Range("AZ5").Value = 100000
$ws.Range("BE5").Value = 110000
$ws.Range("BF5").Value = 250000
$ws.Range("BH5").Value = 250000
$ws.Range("BI5").Value = 250000

# Update sheet view: clear frozen/scrolled topLeftCell and move selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K11").Select()
